$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.83909294746395
$ws.Range("C2").Value = 0.298367494084595
$ws.Range("D2").Value = 26.2732807791768
$ws.Range("E2").Value = 0.0000000000000000000000000000000000000000000000000000000000000000767463425780998

$ws.Range("B3").Value = -0.24309936756535
$ws.Range("C3").Value = 0.604399351405724
$ws.Range("D3").Value = -0.402216460027538
$ws.Range("E3").Value = 0.687981653601289

$ws.Range("B4").Value = -0.166305356584642
$ws.Range("C4").Value = 0.596755938664077
$ws.Range("D4").Value = -0.278682365452349
$ws.Range("E4").Value = 0.780794893625549

$ws.Range("B5").Value = -0.05256878650431
$ws.Range("C5").Value = 1.20646742814453
$ws.Range("D5").Value = -0.0435724871455151
$ws.Range("E5").Value = 0.965291426888526
